# 224614 Add ERP ID update for product sync command
#
# On the "Items" sheet:
#   - mark the existing ITM-1213-3316-0003 row (row 4) for an ERP ID sync
#     update by setting its Action column to "update"
#   - add a new row for ITM-1213-3316-0004, cloned from row 4's data, with
#     a default ("-") Action
#   - extend the Action-column ("C") list validation down to the new row
#   - leave the sheet's selection on C13, matching the saved workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# Flag ITM-1213-3316-0003 (row 4) for an ERP ID update.
$ws.Range("C4").Value = "update"

# Clone row 4 into row 5 (values + styles) for the new item, then set the
# new row's unique ID and reset its Action back to the default.
$ws.Range("A4:R4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = "ITM-1213-3316-0004"
$ws.Range("C5").Value = "-"

# Extend the Action dropdown validation down to the new row.
$ws.Range("C2:C4").Validation.Delete()
$validation = $ws.Range("C2:C5").Validation
$validation.Add(3, 1, 1, '"-,update,review,publish"')
$validation.IgnoreBlank = $false

# Restore the selection Excel left the workbook in when it was saved.
$ws.Range("C13").Select() | Out-Null
